$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.778111702146748
$ws.Range("C2").Value = 0.5310366611625629
$ws.Range("D2").Value = 0.5099897156050792

$ws.Range("B3").Value = 0.7628624588331022
$ws.Range("C3").Value = -0.5589050740759482
$ws.Range("D3").Value = 0.5988320345383111

$ws.Range("B4").Value = -0.8739627011881791
$ws.Range("C4").Value = 0.5053222423535163
$ws.Range("D4").Value = 0.5956262603643074

$ws.Range("B5").Value = 0.779271802461178
$ws.Range("C5").Value = 0.6132461248101666
$ws.Range("D5").Value = -0.6040647386847622

$ws.Range("B6").Value = 0.7739875194658538
$ws.Range("C6").Value = -0.5256880957458819
$ws.Range("D6").Value = 0.6253656874702692

$ws.Range("B7").Value = -0.6244683248714948
$ws.Range("C7").Value = -0.5802676602775768
$ws.Range("D7").Value = 0.5422657080009534

$ws.Range("B8").Value = -0.6910822004340283
$ws.Range("C8").Value = 0.6753835575033513
$ws.Range("D8").Value = 0.5055548986264494

$ws.Range("B9").Value = -0.7012350210805105
$ws.Range("C9").Value = 0.4327636197874499
$ws.Range("D9").Value = 0.5863558691626481
